# Edits the "Hortaliza, Mapocho Venta Directa de Santiago - Pepino dulce" sheet
# Commit: "Fruta / hortaliza, semanal" - adds a new weekly report group
# (date 2022-05-24) as three new rows at the top of the data block
# (rows 101-103), pushing the existing rows 101-209 down to 104-212.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at position 101; this shifts old rows 101:209
# down to 104:212, carrying their values/formatting with them.
$ws.Rows("101:103").Insert()

# Common (constant across this sheet) field values.
$mercadoId = 12
$mercado = "Mapocho Venta Directa de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112043
$categoria = "Pepino dulce"
$variedad = "Cultivar IV Región"
$unidad = "$/bandeja 18 kilos"
$origen = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"
$fecha = 44705   # 2022-05-24 (Excel date serial)

# Row 101: Especial
$ws.Cells.Item(101, 1).Value = $mercadoId
$ws.Cells.Item(101, 2).Value = $mercado
$ws.Cells.Item(101, 3).Value = $region
$ws.Cells.Item(101, 4).Value = $fecha
$ws.Cells.Item(101, 5).Value = $codreg
$ws.Cells.Item(101, 6).Value = $categoriaId
$ws.Cells.Item(101, 7).Value = $categoria
$ws.Cells.Item(101, 8).Value = $variedad
$ws.Cells.Item(101, 9).Value = "Especial"
$ws.Cells.Item(101, 10).Value = 280
$ws.Cells.Item(101, 11).Value = 14000
$ws.Cells.Item(101, 12).Value = 14000
$ws.Cells.Item(101, 13).Value = 14000
$ws.Cells.Item(101, 14).Value = $unidad
$ws.Cells.Item(101, 15).Value = $origen
$ws.Cells.Item(101, 16).Value = 778
$ws.Cells.Item(101, 17).Value = $kgUnidades
$ws.Cells.Item(101, 18).Value = $clasificacion

# Row 102: Primera
$ws.Cells.Item(102, 1).Value = $mercadoId
$ws.Cells.Item(102, 2).Value = $mercado
$ws.Cells.Item(102, 3).Value = $region
$ws.Cells.Item(102, 4).Value = $fecha
$ws.Cells.Item(102, 5).Value = $codreg
$ws.Cells.Item(102, 6).Value = $categoriaId
$ws.Cells.Item(102, 7).Value = $categoria
$ws.Cells.Item(102, 8).Value = $variedad
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 220
$ws.Cells.Item(102, 11).Value = 12000
$ws.Cells.Item(102, 12).Value = 12000
$ws.Cells.Item(102, 13).Value = 12000
$ws.Cells.Item(102, 14).Value = $unidad
$ws.Cells.Item(102, 15).Value = $origen
$ws.Cells.Item(102, 16).Value = 667
$ws.Cells.Item(102, 17).Value = $kgUnidades
$ws.Cells.Item(102, 18).Value = $clasificacion

# Row 103: Segunda
$ws.Cells.Item(103, 1).Value = $mercadoId
$ws.Cells.Item(103, 2).Value = $mercado
$ws.Cells.Item(103, 3).Value = $region
$ws.Cells.Item(103, 4).Value = $fecha
$ws.Cells.Item(103, 5).Value = $codreg
$ws.Cells.Item(103, 6).Value = $categoriaId
$ws.Cells.Item(103, 7).Value = $categoria
$ws.Cells.Item(103, 8).Value = $variedad
$ws.Cells.Item(103, 9).Value = "Segunda"
$ws.Cells.Item(103, 10).Value = 250
$ws.Cells.Item(103, 11).Value = 10000
$ws.Cells.Item(103, 12).Value = 10000
$ws.Cells.Item(103, 13).Value = 10000
$ws.Cells.Item(103, 14).Value = $unidad
$ws.Cells.Item(103, 15).Value = $origen
$ws.Cells.Item(103, 16).Value = 556
$ws.Cells.Item(103, 17).Value = $kgUnidades
$ws.Cells.Item(103, 18).Value = $clasificacion

Write-Host ("New dimension: " + $ws.UsedRange.Address())
